$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be force-written
# as text (NumberFormat "@") and then have their format cleared again so
# the cell keeps matching the original un-styled inline-string cells.

$ws.Range("D2").Value = '27.606.03'
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").Value = '1.597.82'
$ws.Range("E3").Value = '  -1.79%  '
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.63'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("E6").Value = '  -3.45%  '
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("E8").Value = '  -3.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.252'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.60%  '
$ws.Range("E10").Value = '  -3.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("D12").Value = '1.826.61'
$ws.Range("E12").Value = '  -1.53%  '
$ws.Range("D13").Value = '1.591.01'
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.87'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.538'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.55'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.55%  '
$ws.Range("D17").Value = '27.608.04'
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '218.07'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -4.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.41'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.24%  '
$ws.Range("D20").Value = '0.0₃0696'
$ws.Range("E20").Value = '  -3.40%  '
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.20'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.71'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.01'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.52'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.75'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.02%  '
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.10'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.71%  '
$ws.Range("E29").Value = '  -3.60%  '
$ws.Range("E30").Value = '  -1.10%  '
$ws.Range("E31").Value = '  -2.31%  '
$ws.Range("E32").Value = '  -4.09%  '
$ws.Range("D33").Value = '1.373.33'
$ws.Range("E33").Value = '  -1.42%  '
$ws.Range("E34").Value = '  -4.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.972'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.40%  '
$ws.Range("E37").Value = '  -0.75%  '
$ws.Range("E38").Value = '  -2.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.541'
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.816'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.03%  '
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.977'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.38'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("E44").Value = '  -3.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.17'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.32%  '
$ws.Range("D46").Value = '1.736.04'
$ws.Range("E46").Value = '  -1.64%  '
$ws.Range("E47").Value = '  -1.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.95'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").Value = '0.0₆0100'
$ws.Range("E49").Value = '  -3.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0971'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.34%  '
$ws.Range("E51").Value = '  -0.77%  '
